# Applies the "Updated cryptos list" data refresh (Mon May 15 16:12:18 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values that look numeric must be forced to text (leading
# apostrophe) so Excel keeps them as literal strings (e.g. '1.230', matching
# the original inline-string cell content) instead of converting them to
# numbers and losing formatting / significant trailing zeros.

$ws.Range("D2").Value = "27.729.05"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").Value = "1.853.39"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  -1.59%  "

$ws.Range("D5").Value = "'319.58"
$ws.Range("E5").Value = "  -1.56%  "

$ws.Range("E6").Value = "  -1.69%  "

$ws.Range("E7").Value = "  -2.61%  "

$ws.Range("D8").Value = "'0.3749"
$ws.Range("E8").Value = "  -2.18%  "

$ws.Range("D9").Value = "'0.07359"
$ws.Range("E9").Value = "  -1.51%  "

$ws.Range("D10").Value = "'0.8783"
$ws.Range("E10").Value = "  -1.55%  "

$ws.Range("D11").Value = "'21.63"
$ws.Range("E11").Value = "  -0.67%  "

$ws.Range("D12").Value = "1.839.31"
$ws.Range("E12").Value = "  -2.16%  "

$ws.Range("E13").Value = "  -0.83%  "

$ws.Range("D14").Value = "'5.453"
$ws.Range("E14").Value = "  -2.48%  "

$ws.Range("E15").Value = "  -0.77%  "

$ws.Range("D16").Value = "'89.29"
$ws.Range("E16").Value = "  +3.89%  "

$ws.Range("E17").Value = "  -1.89%  "

$ws.Range("D18").Value = "'0.000009006"
$ws.Range("E18").Value = "  -1.78%  "

$ws.Range("D19").Value = "'1.011"
$ws.Range("E19").Value = "  -1.73%  "

$ws.Range("D20").Value = "'15.54"
$ws.Range("E20").Value = "  -0.64%  "

$ws.Range("D21").Value = "27.740.45"
$ws.Range("E21").Value = "  -0.64%  "

$ws.Range("D22").Value = "'5.224"
$ws.Range("E22").Value = "  -2.09%  "

$ws.Range("D23").Value = "'11.10"
$ws.Range("E23").Value = "  -2.26%  "

$ws.Range("D24").Value = "2.080.28"
$ws.Range("E24").Value = "  -1.40%  "

$ws.Range("D25").Value = "'1.995"
$ws.Range("E25").Value = "  -1.59%  "

$ws.Range("D26").Value = "'155.49"
$ws.Range("E26").Value = "  -2.01%  "

$ws.Range("D27").Value = "'18.69"
$ws.Range("E27").Value = "  -1.23%  "

$ws.Range("D28").Value = "'2.207"
$ws.Range("E28").Value = "  +10.61%  "

$ws.Range("D29").Value = "'5.375"

$ws.Range("D30").Value = "'119.29"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").Value = "'0.08942"
$ws.Range("E31").Value = "  -1.40%  "

$ws.Range("D32").Value = "'1.235"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").Value = "'0.7806"
$ws.Range("E33").Value = "  -1.04%  "

$ws.Range("D34").Value = "'4.571"
$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("D35").Value = "'2.926"
$ws.Range("E35").Value = "  -3.44%  "

$ws.Range("E36").Value = "  -1.79%  "

$ws.Range("D37").Value = "'1.135"
$ws.Range("E37").Value = "  -1.32%  "

$ws.Range("D38").Value = "'0.05360"

$ws.Range("D39").Value = "'0.01982"

$ws.Range("D40").Value = "'7.342"
$ws.Range("E40").Value = "  +5.67%  "

$ws.Range("D41").Value = "'2.897"
$ws.Range("E41").Value = "  +0.40%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1698"
$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.5151"
$ws.Range("E43").Value = "  -2.02%  "

$ws.Range("D44").Value = "'8.844"
$ws.Range("E44").Value = "  -1.03%  "

$ws.Range("D45").Value = "'10.79"
$ws.Range("E45").Value = "  +0.10%  "

$ws.Range("D46").Value = "'108.29"
$ws.Range("E46").Value = "  -3.09%  "

$ws.Range("D47").Value = "'0.4798"
$ws.Range("E47").Value = "  +0.64%  "

$ws.Range("D48").Value = "'0.06481"
$ws.Range("E48").Value = "  -2.06%  "

$ws.Range("E49").Value = "  -1.98%  "

$ws.Range("E50").Value = "  -1.90%  "

$ws.Range("D51").Value = "'1.856"
$ws.Range("E51").Value = "  -4.74%  "
